$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  Value = 330000 },
    @{ Row = 3;  Value = 370000 },
    @{ Row = 4;  Value = 120 },
    @{ Row = 5;  Value = 18000 },
    @{ Row = 6;  Value = 690000 },
    @{ Row = 7;  Value = 1600000 },
    @{ Row = 8;  Value = 340000 },
    @{ Row = 9;  Value = 39000 },
    @{ Row = 10; Value = 2200000 },
    @{ Row = 11; Value = 13000 },
    @{ Row = 12; Value = 1000000 },
    @{ Row = 13; Value = 200000 },
    @{ Row = 14; Value = 220000 },
    @{ Row = 15; Value = 190000 },
    @{ Row = 16; Value = 210000 },
    @{ Row = 17; Value = 1600000 },
    @{ Row = 18; Value = 1400000 },
    @{ Row = 19; Value = 7700000 },
    @{ Row = 20; Value = 1200000 },
    @{ Row = 21; Value = 1300000 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.Value
}
